# Add new match-day rows (2025-09-13, Excel serial 45913) to the "Partidos" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partidos")

$startRow = 429
$fecha = 45913

$data = @(
    @("Sebastian Giraldo",            "Amarillo", "Mediocampista", 1, 0, $false, 0, 0, 0, 0, 0),
    @("Carlos Fernando Valencia",     "Amarillo", "Delantero",     2, 0, $false, 0, 0, 0, 0, 0),
    @("Arnul David Narvaez",          "Amarillo", "Delantero",     1, 0, $false, 0, 0, 0, 1, 0),
    @("Oscar Eduardo Herrera",        "Amarillo", "Mediocampista", 1, 0, $false, 0, 1, 0, 0, 0),
    @("Francisco Javier Duran",       "Amarillo", "Defensa",       0, 0, $false, 0, 0, 0, 1, 0),
    @("Edwing Yesid Castillo",        "Amarillo", "Mediocampista", 0, 0, $false, 0, 0, 0, 1, 0),
    @("Andres Jurado",                "Amarillo", "Delantero",     0, 0, $false, 0, 0, 0, 1, 0),
    @("Alexander Uribe",              "Azul",     "Mediocampista", 4, 0, $false, 0, 0, 0, 1, 0),
    @("Andres Tangarife",             "Azul",     "Delantero",     1, 0, $false, 0, 0, 0, 1, 0),
    @("Juan David Espinal",           "Azul",     "Mediocampista", 1, 0, $false, 0, 0, 0, 0, 0),
    @("Jefferson Delgado",            "Azul",     "Mediocampista", 0, 0, $false, 0, 0, 0, 1, 0),
    @("Fabian Grajales",              "Azul",     "Mediocampista", 0, 0, $false, 0, 0, 0, 1, 0),
    @("Andres Guerrero ",             "Azul",     "Defensa",       0, 0, $false, 0, 0, 0, 1, 0),
    @("Fabian Caicedo",               "Amarillo", "Arquero",       0, 0, $true,  6, 0, 0, 0, 0),
    @("Jorge Gonzalez",               "Azul",     "Arquero",       0, 0, $true,  5, 0, 0, 0, 0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value = $fecha
    $ws.Cells.Item($r, 2).Value = $row[0]
    $ws.Cells.Item($r, 3).Value = $row[1]
    $ws.Cells.Item($r, 4).Value = $row[2]
    $ws.Cells.Item($r, 5).Value = $row[3]
    $ws.Cells.Item($r, 6).Value = $row[4]
    $ws.Cells.Item($r, 7).Value = $row[5]
    $ws.Cells.Item($r, 8).Value = $row[6]
    $ws.Cells.Item($r, 9).Value = $row[7]
    $ws.Cells.Item($r, 10).Value = $row[8]
    $ws.Cells.Item($r, 11).Value = $row[9]
    $ws.Cells.Item($r, 12).Value = $row[10]
}

# Update the selected cell on the Partidos sheet view to match the author's final cursor position.
$ws.Range("O418").Select()
